# Refresh the crypto price/volume(1h) columns (D, E) for each coin row
# on Sheet1, per the Mon Sep 30 09:50:07 UTC 2024 GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.850.98'
$ws.Range("E2").Value = '  -2.71%  '
$ws.Range("D3").Value = '2.622.44'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''576.64'
$ws.Range("E5").Value = '  -3.37%  '
$ws.Range("D6").Value = '''155.82'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("D7").Value = '''0.644'
$ws.Range("E7").Value = '  +2.67%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -4.38%  '
$ws.Range("D10").Value = '''5.80'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = '''0.386'
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '''28.47'
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").Value = '3.099.00'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").Value = '''0.0000184'
$ws.Range("E15").Value = '  -6.30%  '
$ws.Range("D16").Value = '63.701.36'
$ws.Range("E16").Value = '  -2.72%  '
$ws.Range("D17").Value = '2.629.50'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '''12.15'
$ws.Range("E18").Value = '  -3.59%  '
$ws.Range("D19").Value = '''4.65'
$ws.Range("E19").Value = '  -1.97%  '
$ws.Range("D20").Value = '''7.55'
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").Value = '''345.21'
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '''67.59'
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("E24").Value = '  +4.96%  '
$ws.Range("E25").Value = '  -2.34%  '
$ws.Range("D26").Value = '''9.31'
$ws.Range("E26").Value = '  -3.90%  '
$ws.Range("D27").Value = '''573.58'
$ws.Range("E27").Value = '  +8.47%  '
$ws.Range("D28").Value = '''1.57'
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("E30").Value = '  -2.11%  '
$ws.Range("D31").Value = '''7.92'
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("D32").Value = '''2.07'
$ws.Range("E32").Value = '  -2.18%  '
$ws.Range("E33").Value = '  -2.33%  '
$ws.Range("D34").Value = '''6.46'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D36").Value = '''0.409'
$ws.Range("E36").Value = '  -2.67%  '
$ws.Range("D37").Value = '''19.94'
$ws.Range("E37").Value = '  -2.12%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").Value = '''151.72'
$ws.Range("E39").Value = '  -2.86%  '
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("D42").Value = '''41.88'
$ws.Range("D43").Value = '''156.48'
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("E44").Value = '  +4.15%  '
$ws.Range("D45").Value = '''3.97'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").Value = '''23.17'
$ws.Range("E46").Value = '  +2.40%  '
$ws.Range("D47").Value = '''0.0596'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("D48").Value = '''0.102'
$ws.Range("E48").Value = '  +2.93%  '
$ws.Range("D49").Value = '''0.631'
$ws.Range("E49").Value = '  -0.48%  '
$ws.Range("D50").Value = '''0.0251'
$ws.Range("E50").Value = '  -1.09%  '
$ws.Range("D51").Value = '''19.12'
$ws.Range("E51").Value = '  -3.01%  '
